# Password.xlsx - "fixed delete password bug"
#
# The sheet layout:
#   Row 1: headers (Domain/Email/Username/Password) + "Counter:" label/value in E1/F1
#          + a merged "flash message" banner H1:R1 used by the password generator macro
#   Row 2: "PrivateKey:" label in E2, derived private key value in F2
#   Row 3: the (now deleted) saved-password entry
#   Row 11 (now 10 after the delete): an anchor cell used by the generator script
#
# The bug: deleting a saved password row left the counter/private-key bookkeeping out of
# sync and left stray formatting behind in the flash-message banner. This change:
#   1. Removes the deleted password's row (row 3), shifting everything below up by one.
#   2. Clears the stray flash-message banner columns (I:R) that used to sit next to the
#      counter cell H1.
#   3. Bumps the Counter value and records the freshly derived PrivateKey.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete the removed password's row.
$ws.Rows(3).Delete()

# 2. Drop the leftover flash-message banner columns next to the counter cell.
$ws.Range("I1:R1").EntireColumn.Delete()

# 3. Refresh the counter and stored private key.
$ws.Range("F1").Value = 4
$ws.Range("F2").Value = "36D7/f4{F5:K1/H8&|85D7/88f4{107F5:85K1/100H8&|113G1@|"
